$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header row had some column labels re-typed in Excel with accidental
# leading/trailing whitespace (this is the actual content change called out
# by the commit message: "Accomodate headers with leading and trailing
# whitespace"). Re-enter the four affected headers with the padding.
$ws.Range("C1").Value = "MIDDLE_NAME  "
$ws.Range("D1").Value = "LAST_NAME  "
$ws.Range("F1").Value = "FAV_NUMBER   "
$ws.Range("G1").Value = "   DATE_REGISTERED"

# Re-apply the (same) body font across every populated cell. This mirrors
# what Excel itself does when the sheet is touched and re-saved: the whole
# stylesheet gets its font explicitly stamped onto every used cell.
for ($r = 1; $r -le 6; $r++) {
    for ($c = 1; $c -le 7; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.Text -ne "") {
            $cell.Font.Name = "Calibri"
        }
    }
}
